$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range("D2").Value = "68.393.37"
$ws.Range("E2").Value = "  -1.35%  "
$ws.Range("D3").Value = "3.823.43"
$ws.Range("E3").Value = "  +1.94%  "
$ws.Range("E4").Value = "  +0.07%  "
Set-TextCell $ws.Range("D5") "599.83"
$ws.Range("E5").Value = "  -0.43%  "
Set-TextCell $ws.Range("D6") "162.97"
$ws.Range("E6").Value = "  -3.28%  "
$ws.Range("D7").Value = "3.825.00"
$ws.Range("E7").Value = "  +2.02%  "
$ws.Range("E8").Value = "  +0.23%  "
Set-TextCell $ws.Range("D9") "0.528"
$ws.Range("E9").Value = "  -2.38%  "
$ws.Range("E10").Value = "  -2.85%  "
Set-TextCell $ws.Range("D11") "6.31"
$ws.Range("E11").Value = "  -1.03%  "
Set-TextCell $ws.Range("D12") "0.458"
$ws.Range("E12").Value = "  -0.84%  "
Set-TextCell $ws.Range("D13") "36.73"
$ws.Range("E13").Value = "  -4.04%  "
Set-TextCell $ws.Range("D14") "0.0000243"
$ws.Range("E14").Value = "  -2.16%  "
$ws.Range("D15").Value = "4.468.89"
$ws.Range("E15").Value = "  +2.07%  "
$ws.Range("D16").Value = "3.806.06"
$ws.Range("E16").Value = "  +1.59%  "
$ws.Range("D17").Value = "68.595.12"
$ws.Range("E17").Value = "  -1.03%  "
Set-TextCell $ws.Range("D18") "7.55"
$ws.Range("E18").Value = "  +1.63%  "
$ws.Range("E19").Value = "  -0.62%  "
Set-TextCell $ws.Range("D20") "17.05"
$ws.Range("E20").Value = "  -1.94%  "
Set-TextCell $ws.Range("D21") "11.17"
$ws.Range("E21").Value = "  -1.23%  "
Set-TextCell $ws.Range("D22") "484.69"
$ws.Range("E22").Value = "  -1.78%  "
Set-TextCell $ws.Range("D23") "0.716"
$ws.Range("E23").Value = "  -1.89%  "
Set-TextCell $ws.Range("D24") "0.0000159"
$ws.Range("E24").Value = "  +6.75%  "
Set-TextCell $ws.Range("D25") "84.00"
$ws.Range("E25").Value = "  -1.04%  "
Set-TextCell $ws.Range("D26") "2.23"
$ws.Range("E26").Value = "  -2.78%  "
Set-TextCell $ws.Range("D27") "12.07"
$ws.Range("E27").Value = "  -2.11%  "
$ws.Range("B28").Value = "Dai"
$ws.Range("C28").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextCell $ws.Range("D28") "0.998"
$ws.Range("E28").Value = "  -0.25%  "
$ws.Range("B29").Value = "RenderToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextCell $ws.Range("D29") "9.98"
$ws.Range("E29").Value = "  -0.97%  "
$ws.Range("E30").Value = "  -1.25%  "
Set-TextCell $ws.Range("D31") "7.82"
$ws.Range("E31").Value = "  -4.48%  "
$ws.Range("D32").Value = "3.980.49"
$ws.Range("E32").Value = "  +2.18%  "
Set-TextCell $ws.Range("D33") "2.35"
$ws.Range("E33").Value = "  -4.46%  "
Set-TextCell $ws.Range("D34") "31.73"
$ws.Range("E34").Value = "  +0.14%  "
$ws.Range("D35").Value = "3.774.47"
$ws.Range("E35").Value = "  +2.46%  "
$ws.Range("E36").Value = "  -1.73%  "
Set-TextCell $ws.Range("D37") "1.02"
$ws.Range("E37").Value = "  +0.97%  "
Set-TextCell $ws.Range("D38") "0.139"
$ws.Range("E38").Value = "  -0.81%  "
Set-TextCell $ws.Range("D39") "5.86"
$ws.Range("E39").Value = "  -2.12%  "
$ws.Range("E40").Value = "  -0.03%  "
$ws.Range("E41").Value = "  -2.91%  "
Set-TextCell $ws.Range("D42") "2.95"
$ws.Range("E42").Value = "  -3.93%  "
Set-TextCell $ws.Range("D43") "427.45"
$ws.Range("E43").Value = "  +0.87%  "
Set-TextCell $ws.Range("D44") "48.47"
$ws.Range("E44").Value = "  -0.94%  "
$ws.Range("E45").Value = "  -0.71%  "
Set-TextCell $ws.Range("D47") "8.40"
$ws.Range("E47").Value = "  -0.88%  "
$ws.Range("D48").Value = "2.838.58"
$ws.Range("E48").Value = "  +1.50%  "
Set-TextCell $ws.Range("D49") "142.67"
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextCell $ws.Range("D50") "25.93"
$ws.Range("E50").Value = "  +12.73%  "
$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextCell $ws.Range("D51") "0.0356"
$ws.Range("E51").Value = "  +0.27%  "
